$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the _GoBack bookmark from the end of the "(2.5 points) ...
#    seeded." paragraph (paragraph 10).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Highlight (green) the "(5 points) ... Postman." paragraph and the
#    "(10 points) ... user experience." paragraph - both currently
#    have no highlighting.
# ------------------------------------------------------------------
$postmanRange = $d.Paragraphs(11).Range
$postmanRange.HighlightColorIndex = 4

$uiRange = $d.Paragraphs(12).Range
$uiRange.HighlightColorIndex = 4

# ------------------------------------------------------------------
# 3. Re-insert the _GoBack bookmark inside the "(5 points) ... see the
#    details of a movie, including title, genre, and director name."
#    paragraph, splitting the run between "genre" and ", and director
#    name."
# ------------------------------------------------------------------
$detailsRange = $d.Paragraphs(13).Range
$findRange = $detailsRange.Duplicate
[void]$findRange.Find.Execute("including title, genre", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$insertPoint = $d.Range($findRange.End, $findRange.End)
[void]$d.Bookmarks.Add("_GoBack", $insertPoint)
